$d = $word.ActiveDocument

# --- Edit 1: merge the split runs of the second
# "The global indicators draw on a wide range..." paragraph (the one
# updated "monthly") into a single run. ---
$monthlyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Figures are updated*monthly*with some data sets being updated more or less frequently*") {
        $monthlyPara = $p
    }
}
$finalText = "The global indicators draw on a wide range of data sets, which are listed at the bottom of the page. Differences between data sets for the same indicator indicate the degree of uncertainty in the indicator. Figures are updated monthly, with some data sets being updated more or less frequently."

$r = $monthlyPara.Range
$r.End = $r.End - 1
$r.Text = "IRON_NATIVE_TMP_PLACEHOLDER"
$r2 = $monthlyPara.Range
$r2.End = $r2.End - 1
$r2.Text = $finalText

# --- Edit 2: insert a new "eei" Heading1 section right before the
# "Terminate" heading (i.e. right after the second
# "precip_quantiles_12month" body paragraph). ---
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*twelve*months aggregated GPCC*") {
        $anchorPara = $p
    }
}

$anchorPara.Range.InsertParagraphAfter()
$headingIndex = $anchorPara.Index + 1
$headingPara = $d.Paragraphs($headingIndex)
$headingPara.Style = "Heading1"
$headingPara.Range.Text = "eei"

$headingPara.Range.InsertParagraphAfter()
$bodyIndex = $headingPara.Index + 1
$bodyPara = $d.Paragraphs($bodyIndex)
$bodyPara.Style = "Normal"
$bodyPara.Range.Text = "Earths energy imbalance is a measure of the net energy flux into the earth system. When the EEI is positive, the amount of energy entering the earth system is larger than the energy leaving the earth system and energy accumulates in the ocean, atmosphere, land and cryosphere, leading to warming. When the EEI is negative, the opposite happens."
